$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the instructions text for quest id 2 ("Exploring the town"), cell D3
$newInstructions = @'
<p>I want you to get to level 10. I do not want you to do this manually. Instead we are going to <a href="/information/automation" target="_blank">explore</a> for an hour. During this time you may gain a Faction level for Surface. The map you are currently on.</p><p><a href="/information/factions" target="_blank">Factions</a> can be seen on your character sheet (tab) under the tab: Factions. As you kill creatures you gain faction points. As you level the faction, you can get what are called <a href="/information/random-enchants" target="_blank">Unique’s</a>. These can be powerful items and you may only have one equipped at a time.</p><p>To do this:</p><p><strong>Desktop:</strong></p><p>- First, find a monster in the drop down list you can kill in one hit.</p><p>- You might want to buy more gear from the shop or better gear, to do this click the hamburger menu at the top left</p><p>- Click Shop, Under General Shop (Uses Gold), click Buy.</p><p>- Here you can buy better gear, eventually you can craft beyond this, how ever right now you might want to spend some of that pretty gold. You can even compare and buy and replace items. How neat!</p><p>- Next, click the Exploration green button (Back on the game page)</p><p>- Select the monster you could kill, 1 hour, Attack.</p><p>- Click explore</p><p><strong>Mobile:</strong></p><p><strong>- </strong>Select the action Fight from the list of actions behind this modal.</p><p>- Find a monster you can kill in one shot.</p><p>- You might want to buy more gear from the shop or better gear, to do this tap the hamburger menu at the top left</p><p>- Tap Shop, Under General Shop (Uses Gold), click Buy.</p><p>- Here you can buy better gear, eventually you can craft beyond this, how ever right now you might want to spend some of that pretty gold. You can even compare and buy and replace items. How neat!</p><p>- Next (Back on the game page), close the fight section, select Exploration from the actions section.</p><p>- Select the monster you could kill, 1 hour, 20 levels and Attack.</p><p>- Tap explore</p><p>Exploring will allow you to log out and idly gain items, exp, gold and possible quest item if the monster drops a quest item. There are things you cannot do while exploring such as change equipment, buy items from the <a href="/information/shop" target="_blank">shop</a> or <a href="/information/market-board" target="_blank">market board</a>. You will be told if you can do an action or not while exploring. Exploring should not just be used to level and log out, you <a href="/information/some-clicking-required" target="_blank">won't get very far</a> if you do that, there is much more to do and we will do it very soon! exploration is a way of gaining levels while doing other things in game.</p><p>Exploration messages will appear, while logged in, in the Exploration chat tab below.</p>
'@
$ws.Range("D3").Value = $newInstructions

# Remove the Labyrinth faction requirement from quest "Go To Labyrinth" (row 8)
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()

# Remove the stale "required_quest_id" references that pointed at quests which no
# longer exist / are no longer required
$ws.Range("T10").ClearContents()
$ws.Range("T12").ClearContents()
$ws.Range("T16").ClearContents()
$ws.Range("T17").ClearContents()
$ws.Range("T18").ClearContents()
$ws.Range("T23").ClearContents()

# Fix the off-by-one id value in the last row (23 -> 22)
$ws.Range("A23").Value = 22

$ws.Columns.Item(20).EntireColumn.AutoFit()
